$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 513-515, pushing the existing 513-520 block
# (two weekly groups of 4 quality rows) down to 516-523.
$ws.Range("A513:A515").EntireRow.Insert()

# Row 513 - new weekly record: Especial
$ws.Cells.Item(513, 1).Value = 9
$ws.Cells.Item(513, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(513, 3).Value = "Metropolitana"
$ws.Cells.Item(513, 4).Value = 45239
$ws.Cells.Item(513, 5).Value = 13
$ws.Cells.Item(513, 6).Value = "Fruta"
$ws.Cells.Item(513, 7).Value = 100107
$ws.Cells.Item(513, 8).Value = "Otros"
$ws.Cells.Item(513, 9).Value = 100107011
$ws.Cells.Item(513, 10).Value = "Tuna"
$ws.Cells.Item(513, 11).Value = "Sin especificar"
$ws.Cells.Item(513, 12).Value = "Especial"
$ws.Cells.Item(513, 13).Value = 200
$ws.Cells.Item(513, 14).Value = 30000
$ws.Cells.Item(513, 15).Value = 30000
$ws.Cells.Item(513, 16).Value = 30000
$ws.Cells.Item(513, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(513, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(513, 19).Value = 1667
$ws.Cells.Item(513, 20).Value = 18

# Row 514 - new weekly record: Primera
$ws.Cells.Item(514, 1).Value = 9
$ws.Cells.Item(514, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(514, 3).Value = "Metropolitana"
$ws.Cells.Item(514, 4).Value = 45239
$ws.Cells.Item(514, 5).Value = 13
$ws.Cells.Item(514, 6).Value = "Fruta"
$ws.Cells.Item(514, 7).Value = 100107
$ws.Cells.Item(514, 8).Value = "Otros"
$ws.Cells.Item(514, 9).Value = 100107011
$ws.Cells.Item(514, 10).Value = "Tuna"
$ws.Cells.Item(514, 11).Value = "Sin especificar"
$ws.Cells.Item(514, 12).Value = "Primera"
$ws.Cells.Item(514, 13).Value = 260
$ws.Cells.Item(514, 14).Value = 25000
$ws.Cells.Item(514, 15).Value = 25000
$ws.Cells.Item(514, 16).Value = 25000
$ws.Cells.Item(514, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(514, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(514, 19).Value = 1389
$ws.Cells.Item(514, 20).Value = 18

# Row 515 - new weekly record: Segunda
$ws.Cells.Item(515, 1).Value = 9
$ws.Cells.Item(515, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(515, 3).Value = "Metropolitana"
$ws.Cells.Item(515, 4).Value = 45239
$ws.Cells.Item(515, 5).Value = 13
$ws.Cells.Item(515, 6).Value = "Fruta"
$ws.Cells.Item(515, 7).Value = 100107
$ws.Cells.Item(515, 8).Value = "Otros"
$ws.Cells.Item(515, 9).Value = 100107011
$ws.Cells.Item(515, 10).Value = "Tuna"
$ws.Cells.Item(515, 11).Value = "Sin especificar"
$ws.Cells.Item(515, 12).Value = "Segunda"
$ws.Cells.Item(515, 13).Value = 220
$ws.Cells.Item(515, 14).Value = 20000
$ws.Cells.Item(515, 15).Value = 20000
$ws.Cells.Item(515, 16).Value = 20000
$ws.Cells.Item(515, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(515, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(515, 19).Value = 1111
$ws.Cells.Item(515, 20).Value = 18
